$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the stray empty cells left on row 11 (C11, D11, F11, G11)
$ws.Range("C11").Value = $null
$ws.Range("D11").Value = $null
$ws.Range("F11").Value = $null
$ws.Range("G11").Value = $null

# Helper: write the date column as literal text rather than letting Excel
# auto-convert the "yyyy-mm-dd" string into a date serial value.
function Set-DateText($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 12
Set-DateText $ws.Range("A12") "2024-05-13"
$ws.Range("B12").Value = "10:47:02"
$ws.Range("C12").Value = "Palet atascado en la curva"

# Row 13
Set-DateText $ws.Range("A13") "2024-05-13"
$ws.Range("B13").Value = "10:47:05"
$ws.Range("D13").Value = "Etiquetadora"

# Row 14
Set-DateText $ws.Range("A14") "2024-05-13"
$ws.Range("B14").Value = "10:47:07"
$ws.Range("F14").Value = "QR desplazado"

# Row 15
Set-DateText $ws.Range("A15") "2024-05-13"
$ws.Range("B15").Value = "10:47:11"
$ws.Range("G15").Value = "Soldadura defectuosa"

# Row 16
Set-DateText $ws.Range("A16") "2024-05-13"
$ws.Range("B16").Value = "10:47:14"
$ws.Range("E16").Value = "No coloca bien el sealling"

# Row 17
Set-DateText $ws.Range("A17") "2024-05-13"
$ws.Range("B17").Value = "10:47:17"
$ws.Range("E17").Value = "Atasco tuerca"

# Row 18 (C18, E18, F18, G18 are left blank)
Set-DateText $ws.Range("A18") "2024-05-13"
$ws.Range("B18").Value = "10:47:56"
$ws.Range("D18").Value = "No coloca bien el sealling"
